$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column G (Тайм код для видео) to match the other columns
$ws.Columns("G").ColumnWidth = 19.17

# Append rows 3-8 with sequential numbering in column A
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7

# Update the selected cell to reflect where the user ended up after entering data
$ws.Range("B6").Select()
